$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the two duplicated header rows (old rows 5 and 3), shifting the
# remaining data rows up so the testcase rows become contiguous.
$ws.Rows.Item(5).Delete() | Out-Null
$ws.Rows.Item(3).Delete() | Out-Null

# Rename the sheet from "Credential" to "Sheet1".
$ws.Name = "Sheet1"

# Move the cell selection from A10 to A11.
$ws.Range("A11").Select() | Out-Null
